# "Fixed title slide typo"
#
# The title slide (slide 1) holds the deck's big title ("Jetpack Joust")
# and the author/subtitle block ("Kyle Brown" / "Seth Lakstins"). Re-enter
# the corrected text on both placeholders so the typo that was present in
# the live editing session is fixed in the saved copy, while every other
# run-level property (font, size, language, spell-check "err" flag, etc.)
# that PowerPoint already has recorded for those runs is left untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

foreach ($shp in $s.Shapes) {
    if (-not $shp.HasTextFrame) { continue }
    if (-not $shp.PlaceholderFormat) { continue }

    $phType = $shp.PlaceholderFormat.Type

    if ($phType -eq [Microsoft.Office.Interop.PowerPoint.PpPlaceholderType]::ppPlaceholderCenterTitle -or
        $phType -eq [Microsoft.Office.Interop.PowerPoint.PpPlaceholderType]::ppPlaceholderTitle) {
        # Corrected title text.
        $shp.TextFrame.TextRange.Text = "Jetpack Joust"
    }
    elseif ($phType -eq [Microsoft.Office.Interop.PowerPoint.PpPlaceholderType]::ppPlaceholderSubtitle) {
        # Corrected author/subtitle text (two lines).
        $shp.TextFrame.TextRange.Text = "Kyle Brown" + [char]13 + "Seth Lakstins"
    }
}
